$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 11 (pushes the existing rows, starting with the
# "004218542" account, down by 2). The two new rows hold updated balances for
# accounts 004550605 (Rejane) and 005142661 (Sabrina), placed in their correct
# descending-Saldo sorted position.
$ws.Rows.Item(11).Insert()
$ws.Rows.Item(11).Insert()

# --- Row 11: 004550605 / Rejane / 16066.57 ---
# Force the account-number cell to be stored as text so the leading zeros are
# kept (plain assignment would otherwise coerce "004550605" to the number
# 4550605). Format as Text, enter the value, then strip the explicit
# NumberFormat stamp back off (via a formats-only paste from an existing
# plain/unstyled text cell) so the cell matches the rest of the "Conta"
# column, which carries no explicit style.
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "004550605"
$ws.Range("A20").Copy()
$ws.Range("A11").PasteSpecial(-4122)

$ws.Range("B11").Value = "Rejane"
$ws.Range("C11").Value = 16066.57

# --- Row 12: 005142661 / Sabrina / 14000.81 ---
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "005142661"
$ws.Range("A20").Copy()
$ws.Range("A12").PasteSpecial(-4122)

$ws.Range("B12").Value = "Sabrina"
$ws.Range("C12").Value = 14000.81

$excel.CutCopyMode = $false

# Remove the two now-stale rows further down the sheet that used to hold
# these accounts' (much smaller) balances. Their original positions were
# row 220 (004550605 / Rejane / 67.8) and row 353 (005142661 / Sabrina /
# 0.81); after inserting 2 rows above them, they now sit at row 222 and
# row 355. Delete bottom-up so the earlier deletion doesn't shift the index
# of the other one.
$ws.Rows.Item(355).Delete()
$ws.Rows.Item(222).Delete()
